$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), formatted like the existing
# header cells (bold font + border + centered), by copying H1's format.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells I2 and J2 in the data row, both value 9.
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
